$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Appliances" worksheet in position 2 (after INFO, before Spring) ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Appliances"

# --- 2. Write header row (row 1) ---
$headers = New-Object 'object[,]' 1,18
$headers[0,0]  = "user"
$headers[0,1]  = "n_users"
$headers[0,2]  = "us_pref"
$headers[0,3]  = "number"
$headers[0,4]  = "P"
$headers[0,5]  = "num_windows"
$headers[0,6]  = "func_time"
$headers[0,7]  = "r_t"
$headers[0,8]  = "func_cycle"
$headers[0,9]  = "fixed"
$headers[0,10] = "fixed_cycle"
$headers[0,11] = "occasional_use"
$headers[0,12] = "flat"
$headers[0,13] = "thermal_P_var"
$headers[0,14] = "pref_index"
$headers[0,15] = "wd_we_type"
$headers[0,16] = "year_min"
$headers[0,17] = "initial_share"
$newSheet.Range("A1:R1").Value = $headers

# --- 3. Write data rows 2-5 ---
$data = New-Object 'object[,]' 4,19
# Row 2
$data[0,0]="0";  $data[0,1]="User1"; $data[0,2]=10; $data[0,3]=1; $data[0,4]=1;  $data[0,5]=100; $data[0,6]=2; $data[0,7]=60;  $data[0,8]=0.1; $data[0,9]=30; $data[0,10]="No";  $data[0,11]=1; $data[0,12]=0.5; $data[0,13]="No";  $data[0,14]=0.1; $data[0,15]=1; $data[0,16]=0; $data[0,17]=2020; $data[0,18]=0.5
# Row 3
$data[1,0]="1";  $data[1,1]="User1"; $data[1,2]=10; $data[1,3]=1; $data[1,4]=2;  $data[1,5]=200; $data[1,6]=3; $data[1,7]=120; $data[1,8]=0.2; $data[1,9]=60; $data[1,10]="Yes"; $data[1,11]=2; $data[1,12]=0.6; $data[1,13]="Yes"; $data[1,14]=0.2; $data[1,15]=2; $data[1,16]=1; $data[1,17]=2021; $data[1,18]=0.6
# Row 4
$data[2,0]="2";  $data[2,1]="User2"; $data[2,2]=20; $data[2,3]=2; $data[2,4]=1;  $data[2,5]=300; $data[2,6]=2; $data[2,7]=60;  $data[2,8]=0.1; $data[2,9]=30; $data[2,10]="No";  $data[2,11]=1; $data[2,12]=0.5; $data[2,13]="No";  $data[2,14]=0.1; $data[2,15]=1; $data[2,16]=0; $data[2,17]=2020; $data[2,18]=0.5
# Row 5
$data[3,0]="3";  $data[3,1]="User2"; $data[3,2]=20; $data[3,3]=2; $data[3,4]=2;  $data[3,5]=400; $data[3,6]=3; $data[3,7]=120; $data[3,8]=0.2; $data[3,9]=60; $data[3,10]="Yes"; $data[3,11]=2; $data[3,12]=0.6; $data[3,13]="Yes"; $data[3,14]=0.2; $data[3,15]=2; $data[3,16]=1; $data[3,17]=2021; $data[3,18]=0.6
$newSheet.Range("A2:S5").Value = $data

Write-Host "data written"

# --- 4. Styling ---
# S1: blank cell with just the dark fill
$newSheet.Range("S1").Interior.Color = 5523012
Write-Host "S1 styled"

# Row1 header A1:R1 - font, fill, alignment common to whole header row
$hdr = $newSheet.Range("A1:R1")
$hdr.Font.Name = "Segoe UI"
$hdr.Font.Size = 9.6
$hdr.Font.Color = 14407121
$hdr.Interior.Color = 5523012
$hdr.HorizontalAlignment = -4108
$hdr.WrapText = $true
Write-Host "header base styled"

# A1:Q1 borders: left, top, bottom = medium D9D9E3 ; no right border
# XlBordersIndex: 7=Left, 8=Top, 9=Bottom, 10=Right
$range_AQ1 = $newSheet.Range("A1:Q1")
$range_AQ1.Borders.Item(7).LineStyle = -4138
$range_AQ1.Borders.Item(7).Weight = -4138
$range_AQ1.Borders.Item(7).Color = 14932441
$range_AQ1.Borders.Item(8).LineStyle = -4138
$range_AQ1.Borders.Item(8).Weight = -4138
$range_AQ1.Borders.Item(8).Color = 14932441
$range_AQ1.Borders.Item(9).LineStyle = -4138
$range_AQ1.Borders.Item(9).Weight = -4138
$range_AQ1.Borders.Item(9).Color = 14932441
Write-Host "A1Q1 borders set"

# R1 borders: all 4 sides medium D9D9E3
$range_R1 = $newSheet.Range("R1")
$range_R1.Borders.Item(7).LineStyle = -4138
$range_R1.Borders.Item(7).Weight = -4138
$range_R1.Borders.Item(7).Color = 14932441
$range_R1.Borders.Item(8).LineStyle = -4138
$range_R1.Borders.Item(8).Weight = -4138
$range_R1.Borders.Item(8).Color = 14932441
$range_R1.Borders.Item(9).LineStyle = -4138
$range_R1.Borders.Item(9).Weight = -4138
$range_R1.Borders.Item(9).Color = 14932441
$range_R1.Borders.Item(10).LineStyle = -4138
$range_R1.Borders.Item(10).Weight = -4138
$range_R1.Borders.Item(10).Color = 14932441
Write-Host "R1 borders set"
